$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: title update
$ws.Range("D9").Value = "국내 대학들 AI, DS 동아리 실태(?)"

# Row 27: title and link update
$ws.Range("D27").Value = "쿠버네티스에서 노드가 추가될 때마다 슬랙 알람 쏘기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/kubernetes-event-alarm/"

# Row 36: title and link update
$ws.Range("D36").Value = "Deep Domain Adaptation"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/374"

# Row 51: title and link update
$ws.Range("D51").Value = "[python] playsound 라이브러리 playsound.PlaysoundException: Error 259 for command 예외 해결 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-playsound-%EB%9D%BC%EC%9D%B4%EB%B8%8C%EB%9F%AC%EB%A6%AC-playsoundPlaysoundException-Error-259-for-command-%EC%98%88%EC%99%B8-%ED%95%B4%EA%B2%B0-%EB%B0%A9%EB%B2%95"
